# Applies a new weekly price record for "Jengibre" (ginger) to the
# "Hortaliza, Mercado Mayorista Lo Valledor de Santiago" subset sheet.
#
# The new record is inserted as row 104 (pushing the previous rows
# 104..141 down to 105..142), which matches how the source data set is
# maintained: newest entries are added near the top of the historical
# block and older rows cascade downward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 104; Excel shifts rows 104:141 down
# to 105:142 automatically and copies row 103's formatting (including the
# date-cell number format) onto the freshly inserted row.
$ws.Rows("104:104").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A104").Value = 6
$ws.Range("B104").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C104").Value = "Metropolitana"
$ws.Range("D104").Value = 45141
$ws.Range("E104").Value = 13
$ws.Range("F104").Value = 100114007
$ws.Range("G104").Value = "Jengibre"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 330
$ws.Range("K104").Value = 16000
$ws.Range("L104").Value = 17000
$ws.Range("M104").Value = 16697
$ws.Range("N104").Value = "$/caja 13 kilos"
$ws.Range("O104").Value = "Perú"
$ws.Range("P104").Value = 1284
$ws.Range("Q104").Value = 13
$ws.Range("R104").Value = "Hortaliza"
